# "Update a few lectures."
#
# 1. Delete the "Key invariant" slide (was slide 8 of 10). Everything after
#    it shifts up by one position, so "The problem" becomes slide 8 and
#    "Another defect" becomes slide 9.
# 2. On slide 2 ("A thread pool"), two paragraphs had their text typed as
#    two separate runs; collapse each back into a single run with the same
#    combined text.

$p = $ppt.ActivePresentation

# --- 1. Remove the "Key invariant" slide -----------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.Item(1).TextFrame.TextRange.Text -eq "Key invariant") {
        $slide.Delete()
        break
    }
}

# --- 2. Merge the split runs on slide 2 -------------------------------------
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2).TextFrame.TextRange

$para1 = $body.Paragraphs(1, 1)
$para1.Text = "TEMP-MERGE-1"
$para1 = $body.Paragraphs(1, 1)
$para1.Text = "A thread pool is a collection of threads which always run."

$para2 = $body.Paragraphs(2, 1)
$para2.Text = "TEMP-MERGE-2"
$para2 = $body.Paragraphs(2, 1)
$para2.Text = "When a task needs to be performed on a thread, one of the threads in the pool will execute it."
